# "got pipeline finally to work and pushed 202912 through"
#
# The genotyping pipeline run completed, so the "genotype" well-metadata
# sheet gets filled in with the (still-to-be-confirmed) genotype call for
# every well, and the workbook is left with that sheet active/selected
# (moving off of "chem_perturbation", which was the sheet left selected
# before).

$wb = $excel.ActiveWorkbook

# --- "genotype" worksheet: fill the plate grid (B2:M9) -------------------
$ws = $wb.Worksheets.Item("genotype")

# All 96 wells get the same provisional call for this run.
$ws.Range("B2:M9").Value = "cep290_unknown"

# Column D needed to be a bit wider to comfortably show the new text.
$ws.Columns.Item(4).ColumnWidth = 17

# Make "genotype" the active sheet with B2 (first well) selected - this is
# where work continues next.
$ws.Activate()
$ws.Range("B2").Select()

# --- "chem_perturbation" worksheet: no longer the active tab -------------
# (Activating "genotype" above already moves tabSelected off of this sheet.)
